# DRS Data is current
# Adds a "Bowler" column (O) to the existing DRS-review log, fixes a
# mis-entered batter name on row 35, and appends seven newly-logged
# DRS reviews (rows 46-52) from matches 15 and 16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Bowler" values for the already-logged reviews (rows 2-45) ---
$ws.Cells.Item(2, 15).Value = "Mustafizur Rahman"
$ws.Cells.Item(6, 15).Value = "AS Joseph"
$ws.Cells.Item(3, 15).Value = "TU Deshpande"
$ws.Cells.Item(4, 15).Value = "TU Deshpande"
$ws.Cells.Item(5, 15).Value = "C Green"
$ws.Cells.Item(7, 15).Value = "HV Patel"
$ws.Cells.Item(8, 15).Value = "Harpreet Brar"
$ws.Cells.Item(9, 15).Value = "Arshdeep Singh"
$ws.Cells.Item(10, 15).Value = "AR Patel"
$ws.Cells.Item(11, 15).Value = "M Markande"
$ws.Cells.Item(12, 15).Value = "Harshit Rana"
$ws.Cells.Item(13, 15).Value = "SP Narine"
$ws.Cells.Item(14, 15).Value = "Ravi Bishnoi"
$ws.Cells.Item(15, 15).Value = "Ravi Bishnoi"
$ws.Cells.Item(16, 15).Value = "Naveen-ul-Haq"
$ws.Cells.Item(17, 15).Value = "Sandeep Sharma"
$ws.Cells.Item(18, 15).Value = "Azmatullah Omarzai"
$ws.Cells.Item(19, 15).Value = "R Sai Kishore"
$ws.Cells.Item(20, 15).Value = "Mohammed Siraj"
$ws.Cells.Item(21, 15).Value = "AS Joseph"
$ws.Cells.Item(22, 15).Value = "AS Joseph"
$ws.Cells.Item(23, 15).Value = "Harpreet Brar"
$ws.Cells.Item(24, 15).Value = "SM Curran"
$ws.Cells.Item(25, 15).Value = "MM Sharma"
$ws.Cells.Item(26, 15).Value = "DL Chahar"
$ws.Cells.Item(27, 15).Value = "TU Deshpande"
$ws.Cells.Item(28, 15).Value = "TU Deshpande"
$ws.Cells.Item(29, 15).Value = "M Pathirana"
$ws.Cells.Item(30, 15).Value = "G Coetzee"
$ws.Cells.Item(31, 15).Value = "SZ Mulani"
$ws.Cells.Item(32, 15).Value = "SZ Mulani"
$ws.Cells.Item(33, 15).Value = "JD Unadkat"

# --- Correction: row 35 batter was mis-entered as "R Ashwin" ---
$ws.Cells.Item(35, 14).Value = "R Parag"

$ws.Cells.Item(34, 15).Value = "Kuldeep Yadav"
$ws.Cells.Item(35, 15).Value = "Kuldeep Yadav"
$ws.Cells.Item(36, 15).Value = "Mukesh Kumar"
$ws.Cells.Item(37, 15).Value = "N Burger"
$ws.Cells.Item(38, 15).Value = "AD Russell"
$ws.Cells.Item(39, 15).Value = "Mohsin Khan"
$ws.Cells.Item(40, 15).Value = "UT Yadav"
$ws.Cells.Item(41, 15).Value = "Noor Ahmad"
$ws.Cells.Item(42, 15).Value = "DG Nalkande"
$ws.Cells.Item(43, 15).Value = "Mustafizur Rahman"
$ws.Cells.Item(44, 15).Value = "TA Boult"
$ws.Cells.Item(45, 15).Value = "R Ashwin"

# --- New DRS reviews (match 15: RCB v LSG, match 16: DC v KKR) ---
# Row 46
$ws.Cells.Item(46, 1).Value = 15
$ws.Cells.Item(46, 2).Value = "RCB"
$ws.Cells.Item(46, 3).Value = "LSG"
$ws.Cells.Item(46, 4).Value = 1
$ws.Cells.Item(46, 5).Value = "LSG"
$ws.Cells.Item(46, 6).Value = "RCB"
$ws.Cells.Item(46, 7).Value = 15
$ws.Cells.Item(46, 8).Value = "RCB"
$ws.Cells.Item(46, 9).Value = "NA Patwardhan"
$ws.Cells.Item(46, 10).Value = "NAP"
$ws.Cells.Item(46, 11).Value = "Wide"
$ws.Cells.Item(46, 12).Value = "Called"
$ws.Cells.Item(46, 13).Value = "Not Called"
$ws.Cells.Item(46, 14).Value = "Q de Kock"
$ws.Cells.Item(46, 15).Value = "Yash Dayal"
$ws.Cells.Item(46, 16).Value = "Successful"
$ws.Cells.Item(46, 17).Value = "No"

# Row 47
$ws.Cells.Item(47, 1).Value = 15
$ws.Cells.Item(47, 2).Value = "RCB"
$ws.Cells.Item(47, 3).Value = "LSG"
$ws.Cells.Item(47, 4).Value = 1
$ws.Cells.Item(47, 5).Value = "LSG"
$ws.Cells.Item(47, 6).Value = "RCB"
$ws.Cells.Item(47, 7).Value = 19
$ws.Cells.Item(47, 8).Value = "LSG"
$ws.Cells.Item(47, 9).Value = "NA Patwardhan"
$ws.Cells.Item(47, 10).Value = "NAP"
$ws.Cells.Item(47, 11).Value = "Wide"
$ws.Cells.Item(47, 12).Value = "Not Called"
$ws.Cells.Item(47, 13).Value = "Not Called"
$ws.Cells.Item(47, 14).Value = "N Pooran"
$ws.Cells.Item(47, 15).Value = "RJW Topley"
$ws.Cells.Item(47, 16).Value = "Unsuccessful"
$ws.Cells.Item(47, 17).Value = "No"

# Row 48
$ws.Cells.Item(48, 1).Value = 15
$ws.Cells.Item(48, 2).Value = "RCB"
$ws.Cells.Item(48, 3).Value = "LSG"
$ws.Cells.Item(48, 4).Value = 1
$ws.Cells.Item(48, 5).Value = "LSG"
$ws.Cells.Item(48, 6).Value = "RCB"
$ws.Cells.Item(48, 7).Value = 20
$ws.Cells.Item(48, 8).Value = "LSG"
$ws.Cells.Item(48, 9).Value = "J Madanagopal"
$ws.Cells.Item(48, 10).Value = "JM"
$ws.Cells.Item(48, 11).Value = "Wide"
$ws.Cells.Item(48, 12).Value = "Not Called"
$ws.Cells.Item(48, 13).Value = "Not Called"
$ws.Cells.Item(48, 14).Value = "N Pooran"
$ws.Cells.Item(48, 15).Value = "Mohammed Siraj"
$ws.Cells.Item(48, 16).Value = "Unsuccessful"
$ws.Cells.Item(48, 17).Value = "No"

# Row 49
$ws.Cells.Item(49, 1).Value = 15
$ws.Cells.Item(49, 2).Value = "RCB"
$ws.Cells.Item(49, 3).Value = "LSG"
$ws.Cells.Item(49, 4).Value = 2
$ws.Cells.Item(49, 5).Value = "RCB"
$ws.Cells.Item(49, 6).Value = "LSG"
$ws.Cells.Item(49, 7).Value = 15
$ws.Cells.Item(49, 8).Value = "RCB"
$ws.Cells.Item(49, 9).Value = "J Madanagopal"
$ws.Cells.Item(49, 10).Value = "JM"
$ws.Cells.Item(49, 11).Value = "Wicket"
$ws.Cells.Item(49, 12).Value = "Out"
$ws.Cells.Item(49, 13).Value = "Not Out"
$ws.Cells.Item(49, 14).Value = "KD Karthik"
$ws.Cells.Item(49, 15).Value = "MP Yadav"
$ws.Cells.Item(49, 16).Value = "Successful"
$ws.Cells.Item(49, 17).Value = "No"

# Row 50
$ws.Cells.Item(50, 1).Value = 15
$ws.Cells.Item(50, 2).Value = "RCB"
$ws.Cells.Item(50, 3).Value = "LSG"
$ws.Cells.Item(50, 4).Value = 2
$ws.Cells.Item(50, 5).Value = "RCB"
$ws.Cells.Item(50, 6).Value = "LSG"
$ws.Cells.Item(50, 7).Value = 19
$ws.Cells.Item(50, 8).Value = "RCB"
$ws.Cells.Item(50, 9).Value = "J Madanagopal"
$ws.Cells.Item(50, 10).Value = "JM"
$ws.Cells.Item(50, 11).Value = "Wicket"
$ws.Cells.Item(50, 12).Value = "Out"
$ws.Cells.Item(50, 13).Value = "Not Out"
$ws.Cells.Item(50, 14).Value = "Mohammad Siraj"
$ws.Cells.Item(50, 15).Value = "Ravi Bishnoi"
$ws.Cells.Item(50, 16).Value = "Successful"
$ws.Cells.Item(50, 17).Value = "No"

# Row 51
$ws.Cells.Item(51, 1).Value = 16
$ws.Cells.Item(51, 2).Value = "DC"
$ws.Cells.Item(51, 3).Value = "KKR"
$ws.Cells.Item(51, 4).Value = 2
$ws.Cells.Item(51, 5).Value = "DC"
$ws.Cells.Item(51, 6).Value = "KKR"
$ws.Cells.Item(51, 7).Value = 2
$ws.Cells.Item(51, 8).Value = "KKR"
$ws.Cells.Item(51, 9).Value = "UV Gandhe"
$ws.Cells.Item(51, 10).Value = "UVG"
$ws.Cells.Item(51, 11).Value = "Wicket"
$ws.Cells.Item(51, 12).Value = "Not Out"
$ws.Cells.Item(51, 13).Value = "Not Out"
$ws.Cells.Item(51, 14).Value = "PP Shaw"
$ws.Cells.Item(51, 15).Value = "VG Arora"
$ws.Cells.Item(51, 16).Value = "Unsuccessful"
$ws.Cells.Item(51, 17).Value = "No"

# Row 52
$ws.Cells.Item(52, 1).Value = 16
$ws.Cells.Item(52, 2).Value = "DC"
$ws.Cells.Item(52, 3).Value = "KKR"
$ws.Cells.Item(52, 4).Value = 2
$ws.Cells.Item(52, 5).Value = "DC"
$ws.Cells.Item(52, 6).Value = "KKR"
$ws.Cells.Item(52, 7).Value = 13
$ws.Cells.Item(52, 8).Value = "KKR"
$ws.Cells.Item(52, 9).Value = "A Totre"
$ws.Cells.Item(52, 10).Value = "AT"
$ws.Cells.Item(52, 11).Value = "Wicket"
$ws.Cells.Item(52, 12).Value = "Not Out"
$ws.Cells.Item(52, 13).Value = "Not Out"
$ws.Cells.Item(52, 14).Value = "Sumit Kumar"
$ws.Cells.Item(52, 15).Value = "CV Varun"
$ws.Cells.Item(52, 16).Value = "Unsuccessful"
$ws.Cells.Item(52, 17).Value = "No"

# --- Widen the new "Bowler" column & refresh the view state ---
$ws.Columns.Item(15).ColumnWidth = 16.3
$ws.Range("M15").Select()
$excel.ActiveWindow.Zoom = 90
